# Append a new trip-log entry as row 11 of the "lorry" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 11

$ws.Cells.Item($row, 1).Value  = "14-03-2018"   # Date
$ws.Cells.Item($row, 2).Value  = "L2"           # vechical name
$ws.Cells.Item($row, 3).Value  = 56000.0        # S.km
$ws.Cells.Item($row, 4).Value  = 57000.0        # C.km
$ws.Cells.Item($row, 5).Value  = 8000.0         # Rent
$ws.Cells.Item($row, 6).Value  = 66.0           # D.rate
$ws.Cells.Item($row, 7).Value  = 150.0          # d.quantity
$ws.Cells.Item($row, 8).Value  = 6.0            # Milage
$ws.Cells.Item($row, 9).Value  = 0.0            # M.rate
$ws.Cells.Item($row, 10).Value = 37600.0        # Total
$ws.Cells.Item($row, 11).Value = ""             # M.Des (blank, like row 2)
